$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39: WAT48 -------------------------------------------------
# Fill values in the same order Excel originally entered them so that
# new shared-string entries land in the expected order:
#   WAT48, description, WAT-323, WAT-558, WAT49, description
$ws.Cells.Item(39, 1).Value = "WAT48"
$ws.Cells.Item(39, 3).Value = "Verify that upon clicking ORCiD tab, Orcid search field should be displayed with an example of orcid number"
$ws.Cells.Item(39, 2).Value = "WAT-323"

# --- Row 40: WAT49 --------------------------------------------------
$ws.Cells.Item(40, 2).Value = "WAT-558"
$ws.Cells.Item(40, 1).Value = "WAT49"
$ws.Cells.Item(40, 3).Value = "Verify that FIND button in ORCid Search page is disabled at the beginning."

# Runmode column for both new rows
$ws.Cells.Item(39, 4).Value = "Y"
$ws.Cells.Item(40, 4).Value = "Y"

# --- Formatting: copy from the existing last row (row 38) so the new
# rows pick up the same borders / wrap-text look used throughout the
# sheet, reusing existing style entries rather than creating new ones.
$ws.Range("A38:E38").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)

$ws.Range("A38:E38").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)

# The description column (C) uses the wrap-text + border style seen on
# other multi-line descriptions in the sheet (e.g. C24).
$ws.Range("C24").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("C40").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View state: matches the selection left behind by the edit -----
$ws.Range("G58").Select()
